# working with cs dataset
# Update the "normal" video feature paths (rows 18-37, column A) from the old
# Colab/Google-Drive path layout to the local WSL path layout, and drop the
# stray placeholder value that had been left in K4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the normal-video feature file paths in column A ---------------
$ws.Range("A18").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_015_x264.npy"
$ws.Range("A19").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_050_x264.npy"
$ws.Range("A20").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_100_x264.npy"
$ws.Range("A21").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_129_x264.npy"
$ws.Range("A22").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_150_x264.npy"
$ws.Range("A23").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_246_x264.npy"
$ws.Range("A24").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_247_x264.npy"
$ws.Range("A25").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_248_x264.npy"
$ws.Range("A26").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_251_x264.npy"
$ws.Range("A27").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_289_x264.npy"
$ws.Range("A28").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_310_x264.npy"
$ws.Range("A29").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_312_x264.npy"
$ws.Range("A30").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_317_x264.npy"
$ws.Range("A31").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_345_x264.npy"
$ws.Range("A32").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_352_x264.npy"
$ws.Range("A33").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_360_x264.npy"
$ws.Range("A34").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_365_x264.npy"
$ws.Range("A35").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_401_x264.npy"
$ws.Range("A36").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_417_x264.npy"
$ws.Range("A37").Value = "/mnt/c/code/data/fea_normal/Normal_Videos_439_x264.npy"

# --- Remove the leftover placeholder value in K4 ---------------------------
$ws.Range("K4").ClearContents()

# --- Let Excel recompute wrapped-text row heights now that several of the
#     paths in column A are shorter than before ----------------------------
$ws.Rows("10:15").AutoFit()
$ws.Rows("18:37").AutoFit()

# --- Update the view: scroll down and select A33 (what the author was
#     looking at when the edit was made) -----------------------------------
$ws.Activate()
$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
